# Actualización automática del inventario: agrega el nuevo producto
# "Almohadilla Epson" (código EZ5LZ1) en la fila 5 de la hoja.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 5

$ws.Cells.Item($row, 1).Value = "EZ5LZ1"
$ws.Cells.Item($row, 2).Value = "Almohadilla Epson"
$ws.Cells.Item($row, 3).Value = "L1110 L1118 L1119 L1210 L1219 L1250 L3100 L3101 L3110 L3115 L3116 L3118 L3150 L3151 L3156 L3158 L3160 L3165 L3166 L3200 L3210 L3218 L3250 L3253 L3260 L3269 L5190 L5290 ET2710 ET2711 ET2712 ET2714 ET2715 ET2720 ET2721 ET2726"
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 100000
$ws.Cells.Item($row, 6).Value = 11
$ws.Cells.Item($row, 7).Value = 23
$ws.Cells.Item($row, 8).Formula = "=(E5-D5)*G5"
$ws.Cells.Item($row, 9).Formula = "=D5*F5"
$ws.Cells.Item($row, 10).Value = 0
